$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 1.02),
    @(2, 3, 1.031462608062972),
    @(2, 4, 1.055817956576718),
    @(2, 5, 1.031031366773808),
    @(2, 6, 1.060440428358929),
    @(2, 9, 1.043140024764572),
    @(2, 10, 1.036598199533664),
    @(2, 11, 1.058557370448738),
    @(2, 12, 1.033840661029319),
    @(2, 13, 1.063167201673782),
    @(2, 14, 1.038070287876477),
    @(3, 2, 1.02),
    @(3, 3, 1.032335070458165),
    @(3, 4, 1.05648050815102),
    @(3, 5, 1.031770875464479),
    @(3, 6, 1.061241999369861),
    @(3, 9, 1.043347729856755),
    @(3, 10, 1.037112939309349),
    @(3, 11, 1.059033823387854),
    @(3, 12, 1.034388833630989),
    @(3, 13, 1.063783231014322),
    @(3, 14, 1.038585758641685),
    @(4, 2, 1.02),
    @(4, 3, 1.032900298392013),
    @(4, 4, 1.056909691339483),
    @(4, 5, 1.032250353067055),
    @(4, 6, 1.061761524108921),
    @(4, 9, 1.04348128003855),
    @(4, 10, 1.037446054578864),
    @(4, 11, 1.059341902763496),
    @(4, 12, 1.034743837368545),
    @(4, 13, 1.064182039974398),
    @(4, 14, 1.038919346973091),
    @(5, 2, 1.02),
    @(5, 3, 1.03313808307667),
    @(5, 4, 1.057090230114572),
    @(5, 5, 1.032452155009907),
    @(5, 6, 1.061980134934369),
    @(5, 9, 1.043537220559029),
    @(5, 10, 1.037586105619529),
    @(5, 11, 1.059471365828225),
    @(5, 12, 1.034893151315929),
    @(5, 13, 1.064349744724485),
    @(5, 14, 1.039059596902298),
    @(6, 2, 1.02),
    @(6, 3, 1.033178017679149),
    @(6, 4, 1.057120549774574),
    @(6, 5, 1.032486051855767),
    @(6, 6, 1.062016852478045),
    @(6, 9, 1.043546601235635),
    @(6, 10, 1.037609621320665),
    @(6, 11, 1.059493100067441),
    @(6, 12, 1.034918225873724),
    @(6, 13, 1.064377905708058),
    @(6, 14, 1.039083145998427),
    @(7, 2, 1.02),
    @(7, 3, 1.032903475044505),
    @(7, 4, 1.056912103276934),
    @(7, 5, 1.032253048654632),
    @(7, 6, 1.061764444403792),
    @(7, 9, 1.04348202832094),
    @(7, 10, 1.037447925912368),
    @(7, 11, 1.05934363286733),
    @(7, 12, 1.034745832233921),
    @(7, 13, 1.064184280676904),
    @(7, 14, 1.038921220964104),
    @(8, 2, 1.02),
    @(8, 3, 1.031757317800408),
    @(8, 4, 1.05604177094676),
    @(8, 5, 1.031281086146764),
    @(8, 6, 1.06071114486792),
    @(8, 9, 1.043210394993907),
    @(8, 10, 1.036772148371763),
    @(8, 11, 1.058718434134835),
    @(8, 12, 1.034025855748448),
    @(8, 13, 1.063375349937837),
    @(8, 14, 1.038244483741878),
    @(9, 2, 1.02),
    @(9, 3, 1.029742964609737),
    @(9, 4, 1.054511801794807),
    @(9, 5, 1.029575838312708),
    @(9, 6, 1.058861730498548),
    @(9, 9, 1.042725273897683),
    @(9, 10, 1.035581733897007),
    @(9, 11, 1.057615151113884),
    @(9, 12, 1.032759515937174),
    @(9, 13, 1.061951490609943),
    @(9, 14, 1.037052378742026),
    @(10, 2, 1.02),
    @(10, 3, 1.028403728346481),
    @(10, 4, 1.053494400583692),
    @(10, 5, 1.02844413376478),
    @(10, 6, 1.057633366757625),
    @(10, 9, 1.042397555456297),
    @(10, 10, 1.034788459454128),
    @(10, 11, 1.056878636504118),
    @(10, 12, 1.031916947572414),
    @(10, 13, 1.061003411676426),
    @(10, 14, 1.036257977758449),
    @(11, 2, 1.02),
    @(11, 3, 1.027824712071377),
    @(11, 4, 1.053054489908936),
    @(11, 5, 1.027955330888932),
    @(11, 6, 1.057102582938013),
    @(11, 9, 1.042254638221914),
    @(11, 10, 1.03444505746281),
    @(11, 11, 1.05655950010274),
    @(11, 12, 1.031552515836742),
    @(11, 13, 1.060593179081401),
    @(11, 14, 1.035914088096914),
    @(12, 2, 1.02),
    @(12, 3, 1.027609773596186),
    @(12, 4, 1.052891184188616),
    @(12, 5, 1.027773954677315),
    @(12, 6, 1.056905594362638),
    @(12, 9, 1.042201400952429),
    @(12, 10, 1.034317517565465),
    @(12, 11, 1.05644092692103),
    @(12, 12, 1.031417211967599),
    @(12, 13, 1.060440845951135),
    @(12, 14, 1.035786367078285),
    @(13, 2, 1.02),
    @(13, 3, 1.027655872543727),
    @(13, 4, 1.05292620941359),
    @(13, 5, 1.027812852002916),
    @(13, 6, 1.056947841438088),
    @(13, 9, 1.042212827368695),
    @(13, 10, 1.034344874613066),
    @(13, 11, 1.056466362682845),
    @(13, 12, 1.031446232266289),
    @(13, 13, 1.060473519837858),
    @(13, 14, 1.035813762976032),
    @(14, 2, 1.02),
    @(14, 3, 1.027806942448964),
    @(14, 4, 1.053040989023546),
    @(14, 5, 1.027940334451072),
    @(14, 6, 1.057086296358563),
    @(14, 9, 1.042250240697735),
    @(14, 10, 1.034434514664317),
    @(14, 11, 1.056549699445112),
    @(14, 12, 1.031541330298108),
    @(14, 13, 1.060580586239026),
    @(14, 14, 1.035903530326438),
    @(15, 2, 1.02),
    @(15, 3, 1.02790003939696),
    @(15, 4, 1.05311172139956),
    @(15, 5, 1.028018905418491),
    @(15, 6, 1.0571716253495),
    @(15, 9, 1.042273272238893),
    @(15, 10, 1.034489746867931),
    @(15, 11, 1.056601041815853),
    @(15, 12, 1.03159993162943),
    @(15, 13, 1.060646559464531),
    @(15, 14, 1.035958840966116),
    @(16, 2, 1.02),
    @(16, 3, 1.028442174454088),
    @(16, 4, 1.053523609460298),
    @(16, 5, 1.028476600121718),
    @(16, 6, 1.057668616603137),
    @(16, 9, 1.0424070191235),
    @(16, 10, 1.034811251935379),
    @(16, 11, 1.056899812009353),
    @(16, 12, 1.031941142371379),
    @(16, 13, 1.06103064375444),
    @(16, 14, 1.036280802607636),
    @(17, 2, 1.02),
    @(17, 3, 1.028782478581629),
    @(17, 4, 1.053782146192618),
    @(17, 5, 1.028764031168922),
    @(17, 6, 1.057980663574458),
    @(17, 9, 1.042490644330353),
    @(17, 10, 1.035012948807122),
    @(17, 11, 1.057087164588343),
    @(17, 12, 1.032155284463462),
    @(17, 13, 1.06127164899873),
    @(17, 14, 1.036482785912072),
    @(18, 2, 1.02),
    @(18, 3, 1.028981057210699),
    @(18, 4, 1.053933007024952),
    @(18, 5, 1.028931803652661),
    @(18, 6, 1.058162781849211),
    @(18, 9, 1.042539323710121),
    @(18, 10, 1.035130603886467),
    @(18, 11, 1.057196422785055),
    @(18, 12, 1.032280229000383),
    @(18, 13, 1.061412251326162),
    @(18, 14, 1.036600608075125),
    @(19, 2, 1.02),
    @(19, 3, 1.029048781762739),
    @(19, 4, 1.053984456912492),
    @(19, 5, 1.028989029834224),
    @(19, 6, 1.058224897486288),
    @(19, 9, 1.042555905511374),
    @(19, 10, 1.035170722668057),
    @(19, 11, 1.057233673334017),
    @(19, 12, 1.032322838446765),
    @(19, 13, 1.061460197792426),
    @(19, 14, 1.036640783829987),
    @(20, 2, 1.02),
    @(20, 3, 1.028745958368002),
    @(20, 4, 1.053754401352474),
    @(20, 5, 1.028733180235386),
    @(20, 6, 1.05794717285136),
    @(20, 9, 1.04248168224644),
    @(20, 10, 1.034991307732344),
    @(20, 11, 1.057067065632028),
    @(20, 12, 1.03213230499211),
    @(20, 13, 1.061245788503858),
    @(20, 14, 1.036461114104486),
    @(21, 2, 1.02),
    @(21, 3, 1.027762452417763),
    @(21, 4, 1.053007186604608),
    @(21, 5, 1.02790278887301),
    @(21, 6, 1.057045520176759),
    @(21, 9, 1.042239227579476),
    @(21, 10, 1.034408117482631),
    @(21, 11, 1.056525159709541),
    @(21, 12, 1.031513324567517),
    @(21, 13, 1.060549056582381),
    @(21, 14, 1.035877095657725),
    @(22, 2, 1.02),
    @(22, 3, 1.027144858911819),
    @(22, 4, 1.052537943413736),
    @(22, 5, 1.027381771309092),
    @(22, 6, 1.056479589331605),
    @(22, 9, 1.042085910861873),
    @(22, 10, 1.034041529386877),
    @(22, 11, 1.056184259184334),
    @(22, 12, 1.03112450813055),
    @(22, 13, 1.060111257770656),
    @(22, 14, 1.035509986964828),
    @(23, 2, 1.02),
    @(23, 3, 1.027472182896457),
    @(23, 4, 1.052786644347615),
    @(23, 5, 1.027657869391437),
    @(23, 6, 1.056779506978738),
    @(23, 9, 1.042167269725618),
    @(23, 10, 1.034235856044593),
    @(23, 11, 1.05636499384345),
    @(23, 12, 1.031330592409967),
    @(23, 13, 1.060343317601579),
    @(23, 14, 1.035704589588687),
    @(24, 2, 1.02),
    @(24, 3, 1.028762460018108),
    @(24, 4, 1.053766937864077),
    @(24, 5, 1.028747120077386),
    @(24, 6, 1.057962305538559),
    @(24, 9, 1.042485732128633),
    @(24, 10, 1.035001086375751),
    @(24, 11, 1.057076147551522),
    @(24, 12, 1.032142688305493),
    @(24, 13, 1.061257473661894),
    @(24, 14, 1.036470906634688),
    @(25, 2, 1.02),
    @(25, 3, 1.030263083479473),
    @(25, 4, 1.054906889281473),
    @(25, 5, 1.030015789492793),
    @(25, 6, 1.059339050032667),
    @(25, 9, 1.042851451282159),
    @(25, 10, 1.035581733897007),
    @(25, 11, 1.057615151113884),
    @(25, 12, 1.032759515937174),
    @(25, 13, 1.061951490609943),
    @(25, 14, 1.037052378742026)
)

foreach ($row in $data) {
    $r = $row[0]
    $c = $row[1]
    $v = $row[2]
    $ws.Cells.Item($r, $c).Value = $v
}
